$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Currency pairs to append below the existing USD/BRL row (row 2).
$pairs = @(
    @("USD", "EUR"),
    @("USD", "JPY"),
    @("USD", "GBP"),
    @("USD", "AUD"),
    @("USD", "CAD"),
    @("USD", "CHF"),
    @("USD", "CNY"),
    @("USD", "SEK"),
    @("USD", "NZD"),
    @("EUR", "BRL"),
    @("EUR", "USD"),
    @("EUR", "GBP"),
    @("EUR", "JPY"),
    @("GBP", "BRL"),
    @("GBP", "USD"),
    @("JPY", "BRL"),
    @("CAD", "BRL"),
    @("AUD", "BRL"),
    @("CHF", "BRL")
)

$row = 3
foreach ($pair in $pairs) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row = $row + 1
}

$ws.Range("A20").Select()
